# Auto-generated COM-interop script replicating the Networking.xlsx data edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture existing cell-format templates (from the current table) into a scratch
#     area before touching any data, so we can re-apply the same named xf entries to
#     the rebuilt rows further down (avoids inventing brand-new/duplicate styles). ---
$ws.Range("A20").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null    # style used by the bold header row
$ws.Range("A21").Copy() | Out-Null
$ws.Range("Z2").PasteSpecial(-4122) | Out-Null    # style used by col-A data cells (bold+border)
$ws.Range("B21").Copy() | Out-Null
$ws.Range("Z3").PasteSpecial(-4122) | Out-Null    # style used by col B:H data cells (border only)
$excel.CutCopyMode = $false

# --- Wipe the old table (rows 1-29) completely: values + formatting ---
$ws.Range("A1:H29").Clear() | Out-Null

# --- Rewrite every row/cell with its final value ---
# Row 1
$ws.Range("A1").Value = "Device"
$ws.Range("B1").Value = "Interface"
$ws.Range("C1").Value = "IP Address "
$ws.Range("D1").Value = "Subnet Mask"
$ws.Range("E1").Value = "Default Gateway"
$ws.Range("F1").Value = "Remarks"

# Row 2
$ws.Range("A2").Value = "NYP-SW-D"
$ws.Range("B2").Value = "VLAN 10"
$ws.Range("C2").Value = "103.31.24.193"
$ws.Range("D2").Value = "255.255.255.224"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "DR1A"

# Row 3
$ws.Range("B3").Value = "VLAN 20"
$ws.Range("C3").Value = "103.31.24.1"
$ws.Range("D3").Value = "255.255.255.192 `t"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "DR1B"

# Row 4
$ws.Range("B4").Value = "VLAN 30"
$ws.Range("C4").Value = "103.31.24.233"
$ws.Range("D4").Value = "255.255.255.248"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "BASEMENT"

# Row 5
$ws.Range("B5").Value = "Gi1/1/1"
$ws.Range("C5").Value = "103.31.24.246"
$ws.Range("D5").Value = "255.255.255.252"

# Row 6
$ws.Range("A6").Value = "NYP-RT"
$ws.Range("B6").Value = "Gi0/0"
$ws.Range("C6").Value = "103.31.24.245"
$ws.Range("D6").Value = "255.255.255.252"

# Row 7
$ws.Range("B7").Value = "Gi0/1"
$ws.Range("C7").Value = "103.31.24.250"
$ws.Range("D7").Value = "255.255.255.252"

# Row 8
$ws.Range("A8").Value = "BASEMENT-DNS"
$ws.Range("B8").Value = "NIC"
$ws.Range("C8").Value = "103.31.24.234"
$ws.Range("D8").Value = "255.255.255.248"
$ws.Range("E8").Value = "103.31.24.233"
$ws.Range("F8").Value = "BASEMENT"

# Row 9
$ws.Range("A9").Value = "BASEMENT-WEB"
$ws.Range("B9").Value = "NIC"
$ws.Range("C9").Value = "103.31.24.235"
$ws.Range("D9").Value = "255.255.255.248"
$ws.Range("E9").Value = "103.31.24.233"
$ws.Range("F9").Value = "BASEMENT"

# Row 10
$ws.Range("A10").Value = "DOVER-SW-D"
$ws.Range("B10").Value = "VLAN 40"
$ws.Range("C10").Value = "103.31.24.65"
$ws.Range("D10").Value = "255.255.255.192"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "SR7A"

# Row 11
$ws.Range("B11").Value = "VLAN 50"
$ws.Range("C11").Value = "103.31.24.129 "
$ws.Range("D11").Value = "255.255.255.192"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "SR7B"

# Row 12
$ws.Range("B12").Value = "VLAN 60"
$ws.Range("C12").Value = "103.31.24.225"
$ws.Range("D12").Value = "255.255.255.248"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "ADMIN"

# Row 13
$ws.Range("B13").Value = "Gi1/1/1"
$ws.Range("C13").Value = "103.31.24.242"
$ws.Range("D13").Value = "255.255.255.252"

# Row 14
$ws.Range("A14").Value = "DOVER-RT"
$ws.Range("B14").Value = "Gi0/0"
$ws.Range("C14").Value = "103.31.24.241"
$ws.Range("D14").Value = "255.255.255.252"

# Row 15
$ws.Range("B15").Value = "Gi0/1"
$ws.Range("C15").Value = "103.31.24.249 "
$ws.Range("D15").Value = "255.255.255.252"

# Row 19
$ws.Range("A19").Value = "Subnet Name"
$ws.Range("B19").Value = "Needed Size"
$ws.Range("C19").Value = "Allocated Size"
$ws.Range("D19").Value = "Address"
$ws.Range("E19").Value = "Mask"
$ws.Range("F19").Value = "Dec Mask"
$ws.Range("G19").Value = "Assignable Range"
$ws.Range("H19").Value = "Broadcast"

# Row 20
$ws.Range("A20").Value = "DR1B"
$ws.Range("B20").Value = 62
$ws.Range("C20").Value = 62
$ws.Range("D20").Value = "103.31.24.0"
$ws.Range("E20").Value = "/26"
$ws.Range("F20").Value = "255.255.255.192"
$ws.Range("G20").Value = "103.31.24.1 - 103.31.24.62"
$ws.Range("H20").Value = "103.31.24.63"

# Row 21
$ws.Range("A21").Value = "SR7A"
$ws.Range("B21").Value = 38
$ws.Range("C21").Value = 62
$ws.Range("D21").Value = "103.31.24.64"
$ws.Range("E21").Value = "/26"
$ws.Range("F21").Value = "255.255.255.192"
$ws.Range("G21").Value = "103.31.24.65 - 103.31.24.126"
$ws.Range("H21").Value = "103.31.24.127"

# Row 22
$ws.Range("A22").Value = "SR7B"
$ws.Range("B22").Value = 38
$ws.Range("C22").Value = 62
$ws.Range("D22").Value = "103.31.24.128"
$ws.Range("E22").Value = "/26"
$ws.Range("F22").Value = "255.255.255.192"
$ws.Range("G22").Value = "103.31.24.129 - 103.31.24.190"
$ws.Range("H22").Value = "103.31.24.191"

# Row 23
$ws.Range("A23").Value = "DR1A"
$ws.Range("B23").Value = 26
$ws.Range("C23").Value = 30
$ws.Range("D23").Value = "103.31.24.192"
$ws.Range("E23").Value = "/27"
$ws.Range("F23").Value = "255.255.255.224"
$ws.Range("G23").Value = "103.31.24.193 - 103.31.24.222"
$ws.Range("H23").Value = "103.31.24.223"

# Row 24
$ws.Range("A24").Value = "ADMIN"
$ws.Range("B24").Value = 6
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = "103.31.24.224"
$ws.Range("E24").Value = "/29"
$ws.Range("F24").Value = "255.255.255.248"
$ws.Range("G24").Value = "103.31.24.225 - 103.31.24.230"
$ws.Range("H24").Value = "103.31.24.231"

# Row 25
$ws.Range("A25").Value = "BASEMENT"
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = "103.31.24.232"
$ws.Range("E25").Value = "/29"
$ws.Range("F25").Value = "255.255.255.248"
$ws.Range("G25").Value = "103.31.24.233 - 103.31.24.238"
$ws.Range("H25").Value = "103.31.24.239"

# Row 26
$ws.Range("A26").Value = "DOVER-RT"
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = "103.31.24.240"
$ws.Range("E26").Value = "/30"
$ws.Range("F26").Value = "255.255.255.252"
$ws.Range("G26").Value = "103.31.24.241 - 103.31.24.242"
$ws.Range("H26").Value = "103.31.24.243"

# Row 27
$ws.Range("A27").Value = "NYP-RT"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = "103.31.24.244"
$ws.Range("E27").Value = "/30"
$ws.Range("F27").Value = "255.255.255.252"
$ws.Range("G27").Value = "103.31.24.245 - 103.31.24.246"
$ws.Range("H27").Value = "103.31.24.247"

# Row 28
$ws.Range("A28").Value = "RT-INTER"
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = "103.31.24.248"
$ws.Range("E28").Value = "/30"
$ws.Range("F28").Value = "255.255.255.252"
$ws.Range("G28").Value = "103.31.24.249 - 103.31.24.250"
$ws.Range("H28").Value = "103.31.24.251"

# --- Apply formatting (re-use the templates captured in column Z) ---

# Row 5-7: new plain-black-font / no-border style (matches the new "Gi x/x" rows)
$ws.Range("B5:D5").Font.Color = 0
$ws.Range("A6:F7").Font.Color = 0

# Row 19: bold header style (same as the template grabbed from the old header row)
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("A19:H19").PasteSpecial(-4122) | Out-Null

# Rows 20-28, column A: bold + border style
$ws.Range("Z2").Copy() | Out-Null
$ws.Range("A20:A28").PasteSpecial(-4122) | Out-Null

# Rows 20-28, columns B-H: border-only style
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B20:H28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row heights: 28.8 for the header, 57.6 for each data row of the subnet table ---
$ws.Range("A19:H19").RowHeight = 28.8
$ws.Range("A20:H28").RowHeight = 57.6

# --- Clean up the scratch template cells ---
$ws.Range("Z1:Z3").Clear() | Out-Null

# --- View state: restore scroll position / zoom / selection similar to the source file ---
$win = $excel.ActiveWindow
$win.Zoom = 130
$win.ScrollRow = 21
$win.ScrollColumn = 1
$ws.Range("B5").Select() | Out-Null

